# Added power options to the arduino side
#
# This edits "Sheet1" of the I2C registries workbook:
#  - fixes a stray trailing space in the "motor 4" label (read table)
#  - adds a new "read killswitch state" read-request row (D8/E8)
#  - adds four new write-request rows (11-14) for desired pitch/roll,
#    default throttle, and enabling/disabling motor power

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix trailing whitespace on the existing "motor 4 " label ---
$ws.Range("E5").Value = "motor 4"

# --- New read-request row: register 57 / "read killswitch state" ---
$ws.Range("D8").Value = 57
$ws.Range("E8").Value = "read killswitch state"

# --- New write-request rows 11-14 ---
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Update desired pitch"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Update desired roll"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Update throttle Default 0"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "enable or disable power to motors"

# --- Match the author's final selection / active cell ---
$ws.Range("B16").Select()

# --- Cosmetic: widen the sheet-tab area (best effort) ---
$wb.Windows.Item(1).TabRatio = 0.976

